# Auto-generated PowerShell COM-interop script
# Applies the cryptos.xlsx price/volume update described in the commit diff
# (GitHub Actions crypto-price refresh, Fri Apr 28 11:02:15 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.292.17'
$ws.Range('E2').Value = '  +1.01%  '
$ws.Range('D3').Value = '1.912.97'
$ws.Range('E3').Value = '  +1.46%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4714'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.81%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4062'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.74'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08045'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.001'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.74'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.01%  '
$ws.Range('D13').Value = '1.911.01'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.886'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.103'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '89.58'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.002'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06624'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = '29.320.52'
$ws.Range('E22').Value = '  +0.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.523'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.45'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.201'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('D26').Value = '2.135.68'
$ws.Range('E26').Value = '  +0.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '153.97'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.81'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.019'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +10.81%  '
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.84'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.073'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09502'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.422'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.96%  '
$ws.Range('E35').Value = '  -1.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.384'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.99%  '
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('E38').Value = '  +1.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.239'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.175'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5854'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.503'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +10.82%  '
$ws.Range('E43').Value = '  +0.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.11'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.07905'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.69%  '
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5506'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.97%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '12.07'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.922'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '113.02'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '44.27'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.32%  '
